# Daily attendance processing - 2026-01-24 06:02:46
# Normalizes the "Recorded By" (column G) values on the
# "Session Analysis Results" sheet: reorders the comma-separated
# list of recorder names for two specific recurring values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "system, backup@backdoor.com, System") {
        $cell.Value = "backup@backdoor.com, system, System"
    }
}
